$wb = $excel.ActiveWorkbook

# --- Step 1: update "总计" (summary) sheet: insert a new row for 2022-Q3 ---
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 16
$summary.Cells.Item(2,4).Value = 1.72
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2

# --- Step 2: create the new "2022-Q3" sheet by copying "2022-Q2" (keeps header/style) ---
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# trim the extra rows carried over from 2022-Q2 (18 data rows) down to the 16 that 2022-Q3 needs
$q3.Range("A18:H19").Delete()

# --- Step 3: overwrite the data cells of "2022-Q3" with the correct values ---
# row 2: rank 0
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).NumberFormat = "@"
$q3.Cells.Item(2,2).Value = "011201"
$q3.Cells.Item(2,3).NumberFormat = "@"
$q3.Cells.Item(2,3).Value = "财通优势行业轮动混合A"
$q3.Cells.Item(2,4).NumberFormat = "@"
$q3.Cells.Item(2,4).Value = "8.53"
$q3.Cells.Item(2,5).NumberFormat = "@"
$q3.Cells.Item(2,5).Value = "91.75"
$q3.Cells.Item(2,6).NumberFormat = "@"
$q3.Cells.Item(2,6).Value = "7.57"
$q3.Cells.Item(2,7).NumberFormat = "@"
$q3.Cells.Item(2,7).Value = "0.6457"
$q3.Cells.Item(2,8).Value = 2

# row 3: rank 1
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).NumberFormat = "@"
$q3.Cells.Item(3,2).Value = "501085"
$q3.Cells.Item(3,3).NumberFormat = "@"
$q3.Cells.Item(3,3).Value = "财通科创主题灵活配置混合（LOF）"
$q3.Cells.Item(3,4).NumberFormat = "@"
$q3.Cells.Item(3,4).Value = "4.83"
$q3.Cells.Item(3,5).NumberFormat = "@"
$q3.Cells.Item(3,5).Value = "86.93"
$q3.Cells.Item(3,6).NumberFormat = "@"
$q3.Cells.Item(3,6).Value = "6.08"
$q3.Cells.Item(3,7).NumberFormat = "@"
$q3.Cells.Item(3,7).Value = "0.2937"
$q3.Cells.Item(3,8).Value = 2

# row 4: rank 2
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).NumberFormat = "@"
$q3.Cells.Item(4,2).Value = "010418"
$q3.Cells.Item(4,3).NumberFormat = "@"
$q3.Cells.Item(4,3).Value = "财通景气行业混合A"
$q3.Cells.Item(4,4).NumberFormat = "@"
$q3.Cells.Item(4,4).Value = "2.72"
$q3.Cells.Item(4,5).NumberFormat = "@"
$q3.Cells.Item(4,5).Value = "94.88"
$q3.Cells.Item(4,6).NumberFormat = "@"
$q3.Cells.Item(4,6).Value = "8.09"
$q3.Cells.Item(4,7).NumberFormat = "@"
$q3.Cells.Item(4,7).Value = "0.2200"
$q3.Cells.Item(4,8).Value = 5

# row 5: rank 3
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).NumberFormat = "@"
$q3.Cells.Item(5,2).Value = "501015"
$q3.Cells.Item(5,3).NumberFormat = "@"
$q3.Cells.Item(5,3).Value = "财通多策略升级混合（LOF）A"
$q3.Cells.Item(5,4).NumberFormat = "@"
$q3.Cells.Item(5,4).Value = "2.06"
$q3.Cells.Item(5,5).NumberFormat = "@"
$q3.Cells.Item(5,5).Value = "94.80"
$q3.Cells.Item(5,6).NumberFormat = "@"
$q3.Cells.Item(5,6).Value = "8.49"
$q3.Cells.Item(5,7).NumberFormat = "@"
$q3.Cells.Item(5,7).Value = "0.1749"
$q3.Cells.Item(5,8).Value = 3

# row 6: rank 4
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).NumberFormat = "@"
$q3.Cells.Item(6,2).Value = "005959"
$q3.Cells.Item(6,3).NumberFormat = "@"
$q3.Cells.Item(6,3).Value = "财通新视野灵活配置混合C"
$q3.Cells.Item(6,4).NumberFormat = "@"
$q3.Cells.Item(6,4).Value = "1.12"
$q3.Cells.Item(6,5).NumberFormat = "@"
$q3.Cells.Item(6,5).Value = "94.59"
$q3.Cells.Item(6,6).NumberFormat = "@"
$q3.Cells.Item(6,6).Value = "8.40"
$q3.Cells.Item(6,7).NumberFormat = "@"
$q3.Cells.Item(6,7).Value = "0.0941"
$q3.Cells.Item(6,8).Value = 5

# row 7: rank 5
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).NumberFormat = "@"
$q3.Cells.Item(7,2).Value = "015271"
$q3.Cells.Item(7,3).NumberFormat = "@"
$q3.Cells.Item(7,3).Value = "财通多策略升级混合（LOF）C"
$q3.Cells.Item(7,4).NumberFormat = "@"
$q3.Cells.Item(7,4).Value = "0.62"
$q3.Cells.Item(7,5).NumberFormat = "@"
$q3.Cells.Item(7,5).Value = "94.80"
$q3.Cells.Item(7,6).NumberFormat = "@"
$q3.Cells.Item(7,6).Value = "8.49"
$q3.Cells.Item(7,7).NumberFormat = "@"
$q3.Cells.Item(7,7).Value = "0.0526"
$q3.Cells.Item(7,8).Value = 3

# row 8: rank 6
$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).NumberFormat = "@"
$q3.Cells.Item(8,2).Value = "005851"
$q3.Cells.Item(8,3).NumberFormat = "@"
$q3.Cells.Item(8,3).Value = "财通新视野灵活配置混合A"
$q3.Cells.Item(8,4).NumberFormat = "@"
$q3.Cells.Item(8,4).Value = "0.62"
$q3.Cells.Item(8,5).NumberFormat = "@"
$q3.Cells.Item(8,5).Value = "94.59"
$q3.Cells.Item(8,6).NumberFormat = "@"
$q3.Cells.Item(8,6).Value = "8.40"
$q3.Cells.Item(8,7).NumberFormat = "@"
$q3.Cells.Item(8,7).Value = "0.0521"
$q3.Cells.Item(8,8).Value = 5

# row 9: rank 7
$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).NumberFormat = "@"
$q3.Cells.Item(9,2).Value = "501032"
$q3.Cells.Item(9,3).NumberFormat = "@"
$q3.Cells.Item(9,3).Value = "财通福盛多策略混合（LOF）A"
$q3.Cells.Item(9,4).NumberFormat = "@"
$q3.Cells.Item(9,4).Value = "0.56"
$q3.Cells.Item(9,5).NumberFormat = "@"
$q3.Cells.Item(9,5).Value = "93.59"
$q3.Cells.Item(9,6).NumberFormat = "@"
$q3.Cells.Item(9,6).Value = "8.03"
$q3.Cells.Item(9,7).NumberFormat = "@"
$q3.Cells.Item(9,7).Value = "0.0450"
$q3.Cells.Item(9,8).Value = 2

# row 10: rank 8
$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).NumberFormat = "@"
$q3.Cells.Item(10,2).Value = "501001"
$q3.Cells.Item(10,3).NumberFormat = "@"
$q3.Cells.Item(10,3).Value = "财通多策略精选混合（LOF）"
$q3.Cells.Item(10,4).NumberFormat = "@"
$q3.Cells.Item(10,4).Value = "0.78"
$q3.Cells.Item(10,5).NumberFormat = "@"
$q3.Cells.Item(10,5).Value = "80.38"
$q3.Cells.Item(10,6).NumberFormat = "@"
$q3.Cells.Item(10,6).Value = "3.72"
$q3.Cells.Item(10,7).NumberFormat = "@"
$q3.Cells.Item(10,7).Value = "0.0290"
$q3.Cells.Item(10,8).Value = 5

# row 11: rank 9
$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).NumberFormat = "@"
$q3.Cells.Item(11,2).Value = "011202"
$q3.Cells.Item(11,3).NumberFormat = "@"
$q3.Cells.Item(11,3).Value = "财通优势行业轮动混合C"
$q3.Cells.Item(11,4).NumberFormat = "@"
$q3.Cells.Item(11,4).Value = "0.34"
$q3.Cells.Item(11,5).NumberFormat = "@"
$q3.Cells.Item(11,5).Value = "91.75"
$q3.Cells.Item(11,6).NumberFormat = "@"
$q3.Cells.Item(11,6).Value = "7.57"
$q3.Cells.Item(11,7).NumberFormat = "@"
$q3.Cells.Item(11,7).Value = "0.0257"
$q3.Cells.Item(11,8).Value = 2

# row 12: rank 10
$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).NumberFormat = "@"
$q3.Cells.Item(12,2).Value = "010637"
$q3.Cells.Item(12,3).NumberFormat = "@"
$q3.Cells.Item(12,3).Value = "财通安盈混合C"
$q3.Cells.Item(12,4).NumberFormat = "@"
$q3.Cells.Item(12,4).Value = "1.55"
$q3.Cells.Item(12,5).NumberFormat = "@"
$q3.Cells.Item(12,5).Value = "48.00"
$q3.Cells.Item(12,6).NumberFormat = "@"
$q3.Cells.Item(12,6).Value = "1.52"
$q3.Cells.Item(12,7).NumberFormat = "@"
$q3.Cells.Item(12,7).Value = "0.0236"
$q3.Cells.Item(12,8).Value = 10

# row 13: rank 11
$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).NumberFormat = "@"
$q3.Cells.Item(13,2).Value = "350009"
$q3.Cells.Item(13,3).NumberFormat = "@"
$q3.Cells.Item(13,3).Value = "天治研究驱动混合A"
$q3.Cells.Item(13,4).NumberFormat = "@"
$q3.Cells.Item(13,4).Value = "0.29"
$q3.Cells.Item(13,5).NumberFormat = "@"
$q3.Cells.Item(13,5).Value = "93.73"
$q3.Cells.Item(13,6).NumberFormat = "@"
$q3.Cells.Item(13,6).Value = "8.09"
$q3.Cells.Item(13,7).NumberFormat = "@"
$q3.Cells.Item(13,7).Value = "0.0235"
$q3.Cells.Item(13,8).Value = 7

# row 14: rank 12
$q3.Cells.Item(14,1).Value = 12
$q3.Cells.Item(14,2).NumberFormat = "@"
$q3.Cells.Item(14,2).Value = "002043"
$q3.Cells.Item(14,3).NumberFormat = "@"
$q3.Cells.Item(14,3).Value = "天治研究驱动混合C"
$q3.Cells.Item(14,4).NumberFormat = "@"
$q3.Cells.Item(14,4).Value = "0.24"
$q3.Cells.Item(14,5).NumberFormat = "@"
$q3.Cells.Item(14,5).Value = "93.73"
$q3.Cells.Item(14,6).NumberFormat = "@"
$q3.Cells.Item(14,6).Value = "8.09"
$q3.Cells.Item(14,7).NumberFormat = "@"
$q3.Cells.Item(14,7).Value = "0.0194"
$q3.Cells.Item(14,8).Value = 7

# row 15: rank 13
$q3.Cells.Item(15,1).Value = 13
$q3.Cells.Item(15,2).NumberFormat = "@"
$q3.Cells.Item(15,2).Value = "010636"
$q3.Cells.Item(15,3).NumberFormat = "@"
$q3.Cells.Item(15,3).Value = "财通安盈混合A"
$q3.Cells.Item(15,4).NumberFormat = "@"
$q3.Cells.Item(15,4).Value = "0.99"
$q3.Cells.Item(15,5).NumberFormat = "@"
$q3.Cells.Item(15,5).Value = "48.00"
$q3.Cells.Item(15,6).NumberFormat = "@"
$q3.Cells.Item(15,6).Value = "1.52"
$q3.Cells.Item(15,7).NumberFormat = "@"
$q3.Cells.Item(15,7).Value = "0.0150"
$q3.Cells.Item(15,8).Value = 10

# row 16: rank 14
$q3.Cells.Item(16,1).Value = 14
$q3.Cells.Item(16,2).NumberFormat = "@"
$q3.Cells.Item(16,2).Value = "014628"
$q3.Cells.Item(16,3).NumberFormat = "@"
$q3.Cells.Item(16,3).Value = "财通福盛多策略混合（LOF）C"
$q3.Cells.Item(16,4).NumberFormat = "@"
$q3.Cells.Item(16,4).Value = "0.09"
$q3.Cells.Item(16,5).NumberFormat = "@"
$q3.Cells.Item(16,5).Value = "93.59"
$q3.Cells.Item(16,6).NumberFormat = "@"
$q3.Cells.Item(16,6).Value = "8.03"
$q3.Cells.Item(16,7).NumberFormat = "@"
$q3.Cells.Item(16,7).Value = "0.0072"
$q3.Cells.Item(16,8).Value = 2

# row 17: rank 15
$q3.Cells.Item(17,1).Value = 15
$q3.Cells.Item(17,2).NumberFormat = "@"
$q3.Cells.Item(17,2).Value = "016234"
$q3.Cells.Item(17,3).NumberFormat = "@"
$q3.Cells.Item(17,3).Value = "财通景气行业混合C"
$q3.Cells.Item(17,4).NumberFormat = "@"
$q3.Cells.Item(17,4).Value = "0.00"
$q3.Cells.Item(17,5).NumberFormat = "@"
$q3.Cells.Item(17,5).Value = "94.88"
$q3.Cells.Item(17,6).NumberFormat = "@"
$q3.Cells.Item(17,6).Value = "8.09"
$q3.Cells.Item(17,7).Value = 0
$q3.Cells.Item(17,8).Value = 5

Write-Output "done"
